$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.478.28'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '1.648.32'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('D4').Formula = "'1.001"
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Formula = "'299.87"
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('D7').Formula = "'0.3794"
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Formula = "'0.3564"
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('D9').Formula = "'50.10"
$ws.Range('E9').Value = '  -3.17%  '
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('D12').Formula = "'1.001"
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Formula = "'22.04"
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').Formula = "'6.410"
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').Formula = "'7.379"
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Formula = "'0.00001197"
$ws.Range('E16').Value = '  -2.95%  '
$ws.Range('D17').Value = '1.653.62'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Formula = "'97.44"
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Formula = "'0.06963"
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Formula = "'6.754"
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').Formula = "'17.28"
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('D22').Formula = "'1.001"
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Formula = "'12.43"
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('D24').Value = '23.504.48'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').Formula = "'2.506"
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Formula = "'2.918"
$ws.Range('E26').Value = '  -5.71%  '
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Formula = "'154.06"
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('D29').Formula = "'5.205"
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Formula = "'132.87"
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').Value = '1.828.10'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Formula = "'6.939"
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').Formula = "'2.124"
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').Formula = "'11.83"
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  -6.54%  '
$ws.Range('D36').Formula = "'0.02728"
$ws.Range('E36').Value = '  -3.33%  '
$ws.Range('D37').Formula = "'0.08737"
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Formula = "'5.963"
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Formula = "'0.2438"
$ws.Range('E39').Value = '  -3.53%  '
$ws.Range('D40').Formula = "'13.19"
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').Formula = "'0.06801"
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('D42').Formula = "'0.6894"
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').Formula = "'1.315"
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('D44').Formula = "'15.48"
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').Formula = "'0.6398"
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('D47').Formula = "'2.260"
$ws.Range('E47').Value = '  -3.78%  '
$ws.Range('D48').Formula = "'3.920"
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('D49').Formula = "'0.07733"
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('D50').Formula = "'127.56"
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').Formula = "'1.152"
$ws.Range('E51').Value = '  -3.56%  '

# Reset style so the forced-text cells do not pick up a quotePrefix-carrying style index
$ws.Range('D4').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
